$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(42613.76190972222, 6, 53, 46, 53, 40, 16382, 16067, 1673, 296, 257, 6, 4, "Bag"),
    @(42613.891898148147, 10, 54, 44, 54, 25, 8474, 14564, 1525, 270, 223, 6, 2, "Bag"),
    @(42614.889270833337, 14, 55, 40, 55, 38, 25322, 24076, 2272, 398, 288, 8, 5, "Bag"),
    @(42615.887824074074, 18, 55, 42, 55, 20, 12136, 18375, 2004, 328, 250, 16, 4, "Bag")
)

$startRow = 8
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    for ($c = 1; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $rowData[$c]
    }
}
